# Scheduled market-data refresh: updates currentAveragePrice* / Leve*Price* / LeveProfit*
# columns (H-N) across several sheets, row by row, per the latest Universalis pull.
$wb = $excel.ActiveWorkbook

# ALC!62
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 33335614
$ws.Range("J62").Value = 2301.5
$ws.Range("L62").Value = 2301.5
$ws.Range("N62").Value = -3549.5

# ALC!65
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 33335614
$ws.Range("J65").Value = 2301.5
$ws.Range("L65").Value = 11507.5
$ws.Range("N65").Value = -17747.5

# ALC!107
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 786.9167
$ws.Range("I107").Value = 756.1111
$ws.Range("J107").Value = 879.3333
$ws.Range("K107").Value = 756.1111
$ws.Range("L107").Value = 879.3333
$ws.Range("M107").Value = 1163.8889
$ws.Range("N107").Value = -4719.3333

# ALC!112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1847.6
$ws.Range("J112").Value = 2554.1428
$ws.Range("L112").Value = 7662.428400000001
$ws.Range("N112").Value = -9878.428400000001

# ALC!113
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 1817.375
$ws.Range("I113").Value = 1453.25
$ws.Range("J113").Value = 1890.2
$ws.Range("K113").Value = 1453.25
$ws.Range("L113").Value = 1890.2
$ws.Range("M113").Value = 1800.75
$ws.Range("N113").Value = -8398.200000000001

# ALC!118
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H118").Value = 670.8333
$ws.Range("I118").Value = 670.8333
$ws.Range("K118").Value = 2012.4999
$ws.Range("M118").Value = -355.4999

# ALC!135
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 376.33334
$ws.Range("I135").Value = 347.5
$ws.Range("J135").Value = 607
$ws.Range("K135").Value = 3127.5
$ws.Range("L135").Value = 5463
$ws.Range("M135").Value = -592.5
$ws.Range("N135").Value = -10533

# ALC!141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 1885.4736
$ws.Range("I141").Value = 1518
$ws.Range("K141").Value = 4554
$ws.Range("M141").Value = 626

# ARM!2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1794.7273
$ws.Range("I2").Value = 1913.4
$ws.Range("J2").Value = 1540.4286
$ws.Range("K2").Value = 1913.4
$ws.Range("L2").Value = 1540.4286
$ws.Range("M2").Value = -1800.4
$ws.Range("N2").Value = -1766.4286

# ARM!32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 29418904
$ws.Range("I32").Value = 38464610
$ws.Range("K32").Value = 38464610
$ws.Range("M32").Value = -38464323

# ARM!45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2457.125
$ws.Range("I45").Value = 2318.0344
$ws.Range("K45").Value = 2318.0344
$ws.Range("M45").Value = -1941.0344

# ARM!110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1223.5714
$ws.Range("I110").Value = 603
$ws.Range("J110").Value = 2775
$ws.Range("K110").Value = 603
$ws.Range("L110").Value = 2775
$ws.Range("M110").Value = 1442
$ws.Range("N110").Value = -6865

# ARM!116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1794.7273
$ws.Range("I116").Value = 1913.4
$ws.Range("J116").Value = 1540.4286
$ws.Range("K116").Value = 1913.4
$ws.Range("L116").Value = 1540.4286
$ws.Range("M116").Value = 380.5999999999999
$ws.Range("N116").Value = -6128.4286

# BSM!3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1794.7273
$ws.Range("I3").Value = 1913.4
$ws.Range("J3").Value = 1540.4286
$ws.Range("K3").Value = 1913.4
$ws.Range("L3").Value = 1540.4286
$ws.Range("M3").Value = -1799.4
$ws.Range("N3").Value = -1768.4286

# CRP!16
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2580.4546
$ws.Range("I16").Value = 1844.7333
$ws.Range("J16").Value = 4157
$ws.Range("K16").Value = 1844.7333
$ws.Range("L16").Value = 4157
$ws.Range("M16").Value = -1557.7333
$ws.Range("N16").Value = -4731

# CRP!99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2261.2
$ws.Range("I99").Value = 1835.3334
$ws.Range("K99").Value = 1835.3334
$ws.Range("M99").Value = -337.3334

# CRP!107
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 522.619
$ws.Range("I107").Value = 357.5
$ws.Range("J107").Value = 742.7778
$ws.Range("K107").Value = 357.5
$ws.Range("L107").Value = 742.7778
$ws.Range("M107").Value = 1562.5
$ws.Range("N107").Value = -4582.7778

# CRP!113
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 2580.4546
$ws.Range("I113").Value = 1844.7333
$ws.Range("J113").Value = 4157
$ws.Range("K113").Value = 1844.7333
$ws.Range("L113").Value = 4157
$ws.Range("M113").Value = 325.2666999999999
$ws.Range("N113").Value = -8497

# CRP!126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 2261.2
$ws.Range("I126").Value = 1835.3334
$ws.Range("K126").Value = 5506.0002
$ws.Range("M126").Value = -3036.0002

# CRP!132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3174.6365
$ws.Range("I132").Value = 806
$ws.Range("J132").Value = 3701
$ws.Range("K132").Value = 2418
$ws.Range("L132").Value = 11103
$ws.Range("M132").Value = 112
$ws.Range("N132").Value = -16163

# CRP!134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 4757.25
$ws.Range("I134").Value = 948.7778
$ws.Range("J134").Value = 9653.857
$ws.Range("K134").Value = 2846.3334
$ws.Range("L134").Value = 28961.571
$ws.Range("M134").Value = -311.3334
$ws.Range("N134").Value = -34031.571

# CUL!80
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 2000
$ws.Range("J80").Value = 2000
$ws.Range("L80").Value = 6000
$ws.Range("N80").Value = -7872

# CUL!83
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H83").Value = 2000
$ws.Range("J83").Value = 2000
$ws.Range("L83").Value = 18000
$ws.Range("N83").Value = -27360

# CUL!131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 744.0769
$ws.Range("I131").Value = 394.5
$ws.Range("J131").Value = 962.5625
$ws.Range("K131").Value = 1183.5
$ws.Range("L131").Value = 2887.6875
$ws.Range("M131").Value = 3856.5
$ws.Range("N131").Value = -12967.6875

# GSM!116
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H116").Value = 42979.332
$ws.Range("J116").Value = 42979.332
$ws.Range("L116").Value = 42979.332
$ws.Range("N116").Value = -52157.332

# LTW!7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3750
$ws.Range("I7").Value = 3750
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 3750
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -3638
$ws.Range("N7").ClearContents()

# LTW!61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 5436.3706
$ws.Range("I61").Value = 6659.1
$ws.Range("K61").Value = 6659.1
$ws.Range("M61").Value = -6457.1

# LTW!113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 5436.3706
$ws.Range("I113").Value = 6659.1
$ws.Range("K113").Value = 6659.1
$ws.Range("M113").Value = -4489.1

# LTW!126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 3750
$ws.Range("I126").Value = 3750
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 11250
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -8780
$ws.Range("N126").ClearContents()

# WVR!113
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 609
$ws.Range("I113").Value = 310
$ws.Range("K113").Value = 930
$ws.Range("M113").Value = 1240

# WVR!138
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H138").Value = 39922.668
$ws.Range("J138").Value = 39922.668
$ws.Range("L138").Value = 39922.668
$ws.Range("N138").Value = -50202.668
